$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New orders pulled in after the "add to cart" bug fix: refresh the order
# rows, extend the table with two more orders, and shift the summary block
# down to make room.

# Keep the order-date column as plain text so values like "2025-06-02" are
# not reinterpreted as serial dates.
$ws.Range("D2:D10").NumberFormat = "@"

$ws.Range("B2").Value = "#ORD#47497504"
$ws.Range("C2").Value = "cez"
$ws.Range("D2").Value = "2025-06-02"
$ws.Range("E2").Value = "Swarm All Black Wireless Gaming Keyboard"
$ws.Range("G2").Value = "₹5,299"
$ws.Range("H2").Value = "₹500"
$ws.Range("B3").Value = "#ORD#33233257"
$ws.Range("C3").Value = "cez"
$ws.Range("D3").Value = "2025-06-02"
$ws.Range("E3").Value = "Hive Full-Size White - Purple Wired Gaming Keyboard"
$ws.Range("G3").Value = "₹3,199"
$ws.Range("H3").Value = "₹100"
$ws.Range("B4").Value = "#ORD#84973799"
$ws.Range("C4").Value = "cez"
$ws.Range("D4").Value = "2025-06-02"
$ws.Range("E4").Value = "Hive Black - Purple Wired Gaming Keyboard"
$ws.Range("G4").Value = "₹2,699"
$ws.Range("H4").Value = "₹100"
$ws.Range("B5").Value = "#ORD#22689663"
$ws.Range("C5").Value = "cez"
$ws.Range("D5").Value = "2025-06-02"
$ws.Range("E5").Value = "Hive Full-Size All White Wired Gaming Keyboard"
$ws.Range("G5").Value = "₹3,149"
$ws.Range("B6").Value = "#ORD#61719044"
$ws.Range("C6").Value = "cez"
$ws.Range("D6").Value = "2025-06-03"
$ws.Range("E6").Value = "Hive Black - Purple Wired Gaming Keyboard"
$ws.Range("G6").Value = "₹2,699"
$ws.Range("B7").Value = "#ORD#61719044"
$ws.Range("C7").Value = "cez"
$ws.Range("D7").Value = "2025-06-03"
$ws.Range("E7").Value = "Hive 75"
$ws.Range("G7").Value = "₹2,799"
$ws.Range("H7").Value = "₹500"
$ws.Range("B8").Value = "#ORD#39594894"
$ws.Range("C8").Value = "cez"
$ws.Range("D8").Value = "2025-06-03"
$ws.Range("E8").Value = "Hive White-Purple Wired Gaming Keyboard GAMING"
$ws.Range("G8").Value = "₹2,699"
$ws.Range("H8").Value = "₹500"
$ws.Range("J8").Value = "Placed"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "#ORD#99299019"
$ws.Range("C9").Value = "cez"
$ws.Range("D9").Value = "2025-06-04"
$ws.Range("E9").Value = "Hive All Black RGB Wired Gaming Keyboard"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = "₹2,649"
$ws.Range("H9").Value = "₹100"
$ws.Range("I9").Value = "razorpay"
$ws.Range("J9").Value = "Delivered"
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "#ORD#99299019"
$ws.Range("C10").Value = "cez"
$ws.Range("D10").Value = "2025-06-04"
$ws.Range("E10").Value = "Hive Black - Purple Wired Gaming Keyboard"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "₹2,699"
$ws.Range("H10").Value = "₹100"
$ws.Range("I10").Value = "razorpay"
$ws.Range("J10").Value = "Delivered"
$ws.Range("A12").Value = "Summary:"
$ws.Range("B12").Value = "Total Orders"
$ws.Range("G12").Value = 10
$ws.Range("B13").Value = "Total Amount (₹)"
$ws.Range("G13").Value = "₹27,891"
$ws.Range("A14").Value = ""
$ws.Range("B14").Value = "Total Discounts (₹)"
$ws.Range("G14").Value = "₹2,900"
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = "Net Sales (₹)"
$ws.Range("G15").Value = "₹24,991"

# Rows that no longer hold data once the summary block moved down
$ws.Range("A11").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("G11").ClearContents()
